# Change "IT-Ausstatung" column values from Ja/Nein to J/N
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blatt1")

# Cells that currently hold "JA" -> set to "J"
$jCells = @("I5","I6","I11","I12","I15","I16","I17","I18","I25","D44","D45","D46","D47","D48","D49")
foreach ($addr in $jCells) {
    $ws.Range($addr).Value = "J"
}

# Cells that currently hold "NEIN" -> set to "N"
$nCells = @("I7","I8","I9","I10","I13","I14","I19","I20","I21","I22","I23","I24","I26","I27","I28","I29","I30","I31","I32","I33","I34","I35","D50","D51","D52","D53","D54","D55")
foreach ($addr in $nCells) {
    $ws.Range($addr).Value = "N"
}
